$d = $word.ActiveDocument

$d.Content.Find.Execute(" uma classe de associação para ", $true, $false, $false, $false, $false, $true, 1, $false, " uma classe para ", 2)
